$wb = $excel.ActiveWorkbook

# Sheet: Login -> G2/G3 "Final Status" column
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("G2").Value = "Success - 2020/12/19 16:28:07"
$wsLogin.Range("G3").Value = "Success - 2020/12/19 16:28:11"

# Sheet: School Search -> C2/C3 status column
$wsSchool = $wb.Worksheets.Item("School Search")
$wsSchool.Range("C2").Value = "Success - 2020/12/19 16:28:18"
$wsSchool.Range("C3").Value = "Success - 2020/12/19 16:28:21"

# Sheet: Product Search -> K2/K3/K4 "Final Result" column
$wsProduct = $wb.Worksheets.Item("Product Search")
$wsProduct.Range("K2").Value = "Success - 2020/12/19 16:28:48"
$wsProduct.Range("K3").Value = "Success - 2020/12/19 16:29:10"
$wsProduct.Range("K4").Value = "Success - 2020/12/19 16:29:33"

# Sheet: Shopping Cart -> G2/G3/G4 "Final Result" column
$wsCart = $wb.Worksheets.Item("Shopping Cart")
$wsCart.Range("G2").Value = "Success - 2020/12/19 16:29:35"
$wsCart.Range("G3").Value = "Success - 2020/12/19 16:29:35"
$wsCart.Range("G4").Value = "Success - 2020/12/19 16:29:36"

# Sheet: Checkout -> P2/P3/P4 status column
$wsCheckout = $wb.Worksheets.Item("Checkout")
$wsCheckout.Range("P2").Value = "Success - 2020/12/19 16:29:46"
$wsCheckout.Range("P3").Value = "Success - 2020/12/19 16:29:53"
$wsCheckout.Range("P4").Value = "Success - 2020/12/19 16:30:02"

# Sheet: Payment -> C2 status column (cell uses a quote-prefixed / text-forced
# style, so assign via Formula with a leading apostrophe to keep that format)
$wsPayment = $wb.Worksheets.Item("Payment")
$wsPayment.Range("C2").Formula = "'Success - 2020/12/19 16:30:12"
